$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.396.96'
$ws.Range("E2").Value = '  -0.74%  '

$ws.Range("D3").Value = '1.861.72'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.009'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.56%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.59'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.42%  '

$ws.Range("E7").Value = '  -1.97%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3826'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07811'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.22%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9842'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.48'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.59%  '

$ws.Range("D12").Value = '1.889.01'
$ws.Range("E12").Value = '  +3.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.883'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.98%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.626'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.94%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06905'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.009'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.56%  '

$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.68'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.66%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009923'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.63'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.29%  '

$ws.Range("D21").Value = '28.415.51'
$ws.Range("E21").Value = '  -0.56%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.240'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.85'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.099'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.56%  '

$ws.Range("D25").Value = '2.093.06'
$ws.Range("E25").Value = '  +2.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.645'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.31'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.83%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.910'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.40%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09271'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.51%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9034'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.15%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.255'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.314'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.83%  '

$ws.Range("E35").Value = '  -1.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05682'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.144'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02054'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.64%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.641'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.50%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5549'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.86%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1764'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.599'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07089'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.49'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5227'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.127'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.106'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.805'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.41%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '111.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.423'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.007'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
